$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 141
$ws.Range("D141").Value = 44524
$ws.Range("J141").Value = 200
$ws.Range("K141").Value = 18000
$ws.Range("L141").Value = 20000
$ws.Range("M141").Value = 19000
$ws.Range("P141").Value = 760

# Row 142
$ws.Range("D142").Value = 44524
$ws.Range("I142").Value = 'Segunda'
$ws.Range("J142").Value = 120
$ws.Range("K142").Value = 10000
$ws.Range("L142").Value = 12000
$ws.Range("M142").Value = 11000
$ws.Range("P142").Value = 440

# Row 143
$ws.Range("D143").Value = 44524
$ws.Range("H143").Value = 'Inferno'
$ws.Range("I143").Value = 'Primera'
$ws.Range("J143").Value = 160
$ws.Range("K143").Value = 33000
$ws.Range("L143").Value = 35000
$ws.Range("M143").Value = 34000
$ws.Range("P143").Value = 1360

# Row 144
$ws.Range("D144").Value = 44426
$ws.Range("J144").Value = 240
$ws.Range("K144").Value = 65000
$ws.Range("L144").Value = 70000
$ws.Range("M144").Value = 67500
$ws.Range("P144").Value = 2700

# Row 145
$ws.Range("D145").Value = 44448
$ws.Range("J145").Value = 160
$ws.Range("K145").Value = 65000
$ws.Range("L145").Value = 70000
$ws.Range("M145").Value = 67500
$ws.Range("P145").Value = 2700

# Row 146
$ws.Range("D146").Value = 44448
$ws.Range("H146").Value = 'Americana (o)'
$ws.Range("I146").Value = 'Segunda'
$ws.Range("J146").Value = 100
$ws.Range("K146").Value = 55000
$ws.Range("L146").Value = 60000
$ws.Range("M146").Value = 57500
$ws.Range("P146").Value = 2300

# Row 147
$ws.Range("D147").Value = 44210
$ws.Range("J147").Value = 100
$ws.Range("K147").Value = 17000
$ws.Range("L147").Value = 20000
$ws.Range("M147").Value = 18500
$ws.Range("P147").Value = 740

# Row 148
$ws.Range("D148").Value = 44301
$ws.Range("I148").Value = 'Primera'
$ws.Range("J148").Value = 300
$ws.Range("K148").Value = 15000
$ws.Range("L148").Value = 17000
$ws.Range("M148").Value = 16000
$ws.Range("P148").Value = 640

# Row 149
$ws.Range("D149").Value = 44301
$ws.Range("H149").Value = 'Inferno'
$ws.Range("J149").Value = 160
$ws.Range("K149").Value = 18000
$ws.Range("L149").Value = 19000
$ws.Range("M149").Value = 18500
$ws.Range("P149").Value = 740

# Row 150
$ws.Range("D150").Value = 44441
$ws.Range("H150").Value = 'Americana (o)'
$ws.Range("J150").Value = 160
$ws.Range("K150").Value = 65000
$ws.Range("L150").Value = 70000
$ws.Range("M150").Value = 67500
$ws.Range("P150").Value = 2700

# Row 151
$ws.Range("D151").Value = 44441
$ws.Range("I151").Value = 'Segunda'
$ws.Range("J151").Value = 100
$ws.Range("K151").Value = 55000
$ws.Range("L151").Value = 60000
$ws.Range("M151").Value = 57500
$ws.Range("P151").Value = 2300

# Row 152
$ws.Range("D152").Value = 44279
$ws.Range("I152").Value = 'Primera'
$ws.Range("J152").Value = 240
$ws.Range("K152").Value = 11000
$ws.Range("L152").Value = 12000
$ws.Range("M152").Value = 11500
$ws.Range("P152").Value = 460

# Row 153
$ws.Range("D153").Value = 44279
$ws.Range("H153").Value = 'Inferno'
$ws.Range("J153").Value = 200
$ws.Range("K153").Value = 17000
$ws.Range("L153").Value = 18000
$ws.Range("M153").Value = 17500
$ws.Range("P153").Value = 700

# Row 154
$ws.Range("D154").Value = 44504
$ws.Range("H154").Value = 'Americana (o)'
$ws.Range("J154").Value = 240
$ws.Range("K154").Value = 30000
$ws.Range("L154").Value = 35000
$ws.Range("M154").Value = 32500
$ws.Range("P154").Value = 1300

# Row 155
$ws.Range("D155").Value = 44504
$ws.Range("I155").Value = 'Segunda'
$ws.Range("J155").Value = 100
$ws.Range("K155").Value = 20000
$ws.Range("L155").Value = 25000
$ws.Range("M155").Value = 22500
$ws.Range("P155").Value = 900

# Row 156
$ws.Range("D156").Value = 44350
$ws.Range("H156").Value = 'Americana (o)'
$ws.Range("J156").Value = 240
$ws.Range("K156").Value = 33000
$ws.Range("L156").Value = 35000
$ws.Range("M156").Value = 34000
$ws.Range("P156").Value = 1360

# Row 157
$ws.Range("D157").Value = 44350
$ws.Range("I157").Value = 'Primera'
$ws.Range("J157").Value = 160
$ws.Range("L157").Value = 37000
$ws.Range("M157").Value = 36000
$ws.Range("P157").Value = 1440

# Row 158
$ws.Range("D158").Value = 44384
$ws.Range("J158").Value = 140
$ws.Range("K158").Value = 45000
$ws.Range("L158").Value = 50000
$ws.Range("M158").Value = 47500
$ws.Range("P158").Value = 1900

# Row 159
$ws.Range("D159").Value = 44384
$ws.Range("J159").Value = 80
$ws.Range("K159").Value = 48000
$ws.Range("L159").Value = 50000
$ws.Range("M159").Value = 49000
$ws.Range("P159").Value = 1960

# Row 160
$ws.Range("D160").Value = 44384
$ws.Range("H160").Value = 'Inferno'
$ws.Range("I160").Value = 'Segunda'
$ws.Range("J160").Value = 30
$ws.Range("K160").Value = 35000
$ws.Range("L160").Value = 40000
$ws.Range("M160").Value = 37500
$ws.Range("P160").Value = 1500

# Row 161
$ws.Range("D161").Value = 44329
$ws.Range("J161").Value = 160
$ws.Range("K161").Value = 22000
$ws.Range("L161").Value = 24000
$ws.Range("M161").Value = 23000
$ws.Range("P161").Value = 920

# Row 162
$ws.Range("D162").Value = 44329
$ws.Range("H162").Value = 'Inferno'
$ws.Range("I162").Value = 'Primera'
$ws.Range("J162").Value = 160
$ws.Range("K162").Value = 28000
$ws.Range("L162").Value = 30000
$ws.Range("M162").Value = 29000
$ws.Range("P162").Value = 1160

# Row 163
$ws.Range("D163").Value = 44272
$ws.Range("K163").Value = 10000
$ws.Range("L163").Value = 12000
$ws.Range("M163").Value = 11000
$ws.Range("P163").Value = 440

# Row 164
$ws.Range("D164").Value = 44447
$ws.Range("H164").Value = 'Americana (o)'
$ws.Range("J164").Value = 140
$ws.Range("K164").Value = 65000
$ws.Range("L164").Value = 70000
$ws.Range("M164").Value = 67500
$ws.Range("P164").Value = 2700

# Row 165
$ws.Range("D165").Value = 44447
$ws.Range("I165").Value = 'Segunda'
$ws.Range("J165").Value = 100
$ws.Range("K165").Value = 55000
$ws.Range("L165").Value = 60000
$ws.Range("M165").Value = 57500
$ws.Range("P165").Value = 2300

# Row 166
$ws.Range("D166").Value = 44315
$ws.Range("J166").Value = 300

# Row 167
$ws.Range("D167").Value = 44315
$ws.Range("J167").Value = 240
$ws.Range("K167").Value = 23000
$ws.Range("L167").Value = 25000
$ws.Range("M167").Value = 24000
$ws.Range("P167").Value = 960

# Row 168
$ws.Range("D168").Value = 44322
$ws.Range("K168").Value = 23000
$ws.Range("L168").Value = 25000
$ws.Range("M168").Value = 24000
$ws.Range("P168").Value = 960

# Row 169
$ws.Range("D169").Value = 44510
$ws.Range("J169").Value = 240
$ws.Range("K169").Value = 20000
$ws.Range("L169").Value = 22000
$ws.Range("M169").Value = 21000
$ws.Range("P169").Value = 840

# Row 170
$ws.Range("D170").Value = 44510
$ws.Range("K170").Value = 33000
$ws.Range("L170").Value = 35000
$ws.Range("M170").Value = 34000
$ws.Range("P170").Value = 1360

# Row 171
$ws.Range("D171").Value = 44161
$ws.Range("J171").Value = 200
$ws.Range("K171").Value = 17000
$ws.Range("L171").Value = 18000
$ws.Range("M171").Value = 17500
$ws.Range("P171").Value = 700

# Row 172
$ws.Range("D172").Value = 44468
$ws.Range("I172").Value = 'Primera'
$ws.Range("K172").Value = 80000
$ws.Range("L172").Value = 85000
$ws.Range("M172").Value = 82500
$ws.Range("P172").Value = 3300

# Row 173
$ws.Range("D173").Value = 44468
$ws.Range("J173").Value = 160
$ws.Range("K173").Value = 65000
$ws.Range("L173").Value = 70000
$ws.Range("M173").Value = 67500
$ws.Range("P173").Value = 2700

# Row 174
$ws.Range("D174").Value = 44517
$ws.Range("J174").Value = 400
$ws.Range("K174").Value = 20000
$ws.Range("L174").Value = 22000
$ws.Range("M174").Value = 21000
$ws.Range("P174").Value = 840

# Row 175
$ws.Range("D175").Value = 44517
$ws.Range("H175").Value = 'Americana (o)'
$ws.Range("I175").Value = 'Segunda'
$ws.Range("J175").Value = 200
$ws.Range("K175").Value = 12000
$ws.Range("L175").Value = 15000
$ws.Range("M175").Value = 13500
$ws.Range("P175").Value = 540

# Row 176
$ws.Range("D176").Value = 44517
$ws.Range("H176").Value = 'Inferno'
$ws.Range("J176").Value = 240
$ws.Range("K176").Value = 28000
$ws.Range("L176").Value = 30000
$ws.Range("M176").Value = 29000
$ws.Range("P176").Value = 1160

# Row 177
$ws.Range("D177").Value = 44238
$ws.Range("I177").Value = 'Primera'
$ws.Range("J177").Value = 160
$ws.Range("K177").Value = 7000
$ws.Range("L177").Value = 8000
$ws.Range("M177").Value = 7500
$ws.Range("P177").Value = 300

# Row 178
$ws.Range("D178").Value = 44238
$ws.Range("H178").Value = 'Inferno'
$ws.Range("J178").Value = 160
$ws.Range("K178").Value = 16000
$ws.Range("L178").Value = 18000
$ws.Range("M178").Value = 17000
$ws.Range("P178").Value = 680

# Row 179
$ws.Range("D179").Value = 44391
$ws.Range("H179").Value = 'Americana (o)'
$ws.Range("K179").Value = 60000
$ws.Range("L179").Value = 65000
$ws.Range("M179").Value = 62500
$ws.Range("P179").Value = 2500

# Row 180
$ws.Range("D180").Value = 44391
$ws.Range("I180").Value = 'Segunda'
$ws.Range("J180").Value = 60
$ws.Range("K180").Value = 50000
$ws.Range("L180").Value = 55000
$ws.Range("M180").Value = 52500
$ws.Range("P180").Value = 2100

# Row 181
$ws.Range("D181").Value = 44251
$ws.Range("H181").Value = 'Americana (o)'
$ws.Range("K181").Value = 6000
$ws.Range("L181").Value = 8000
$ws.Range("M181").Value = 7000
$ws.Range("P181").Value = 280

# Row 182
$ws.Range("D182").Value = 44251
$ws.Range("H182").Value = 'Inferno'
$ws.Range("J182").Value = 160
$ws.Range("K182").Value = 12000
$ws.Range("L182").Value = 14000
$ws.Range("M182").Value = 13000
$ws.Range("P182").Value = 520

# Row 183
$ws.Range("D183").Value = 44181
$ws.Range("H183").Value = 'Americana (o)'
$ws.Range("J183").Value = 300
$ws.Range("K183").Value = 16000
$ws.Range("L183").Value = 21000
$ws.Range("M183").Value = 18500
$ws.Range("P183").Value = 740

# Row 184 (new)
$ws.Range("D184").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("A184").Value = 2
$ws.Range("B184").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C184").Value = 'Coquimbo'
$ws.Range("D184").Value = 44181
$ws.Range("E184").Value = 4
$ws.Range("F184").Value = 100112021
$ws.Range("G184").Value = 'Ají'
$ws.Range("H184").Value = 'Inferno'
$ws.Range("I184").Value = 'Primera'
$ws.Range("J184").Value = 200
$ws.Range("K184").Value = 25000
$ws.Range("L184").Value = 27000
$ws.Range("M184").Value = 26000
$ws.Range("N184").Value = '$/caja 25 kilos'
$ws.Range("O184").Value = 'Provincia de Limarí'
$ws.Range("P184").Value = 1040
$ws.Range("Q184").Value = 25
$ws.Range("R184").Value = 'Hortaliza'

# Row 185 (new)
$ws.Range("D185").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("A185").Value = 2
$ws.Range("B185").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C185").Value = 'Coquimbo'
$ws.Range("D185").Value = 44307
$ws.Range("E185").Value = 4
$ws.Range("F185").Value = 100112021
$ws.Range("G185").Value = 'Ají'
$ws.Range("H185").Value = 'Americana (o)'
$ws.Range("I185").Value = 'Primera'
$ws.Range("J185").Value = 200
$ws.Range("K185").Value = 16000
$ws.Range("L185").Value = 18000
$ws.Range("M185").Value = 17000
$ws.Range("N185").Value = '$/caja 25 kilos'
$ws.Range("O185").Value = 'Provincia de Limarí'
$ws.Range("P185").Value = 680
$ws.Range("Q185").Value = 25
$ws.Range("R185").Value = 'Hortaliza'

# Row 186 (new)
$ws.Range("D186").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("A186").Value = 2
$ws.Range("B186").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C186").Value = 'Coquimbo'
$ws.Range("D186").Value = 44307
$ws.Range("E186").Value = 4
$ws.Range("F186").Value = 100112021
$ws.Range("G186").Value = 'Ají'
$ws.Range("H186").Value = 'Inferno'
$ws.Range("I186").Value = 'Primera'
$ws.Range("J186").Value = 200
$ws.Range("K186").Value = 23000
$ws.Range("L186").Value = 25000
$ws.Range("M186").Value = 24000
$ws.Range("N186").Value = '$/caja 25 kilos'
$ws.Range("O186").Value = 'Provincia de Limarí'
$ws.Range("P186").Value = 960
$ws.Range("Q186").Value = 25
$ws.Range("R186").Value = 'Hortaliza'

